$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 932.1667
$ws.Range("I98").Value = 932.1667
$ws.Range("K98").Value = 932.1667
$ws.Range("M98").Value = 565.8333
$ws.Range("H122").Value = 932.1667
$ws.Range("I122").Value = 932.1667
$ws.Range("K122").Value = 2796.5001
$ws.Range("M122").Value = -346.5001000000002
$ws.Range("H138").Value = 2613.918
$ws.Range("J138").Value = 3430
$ws.Range("L138").Value = 10290
$ws.Range("N138").Value = -20570
$ws.Range("H141").Value = 4172.7646
$ws.Range("I141").Value = 4130.533
$ws.Range("J141").Value = 4489.5
$ws.Range("K141").Value = 12391.599
$ws.Range("L141").Value = 13468.5
$ws.Range("M141").Value = -7211.599000000002
$ws.Range("N141").Value = -23828.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 15800.8
$ws.Range("J2").Value = 1006
$ws.Range("L2").Value = 1006
$ws.Range("N2").Value = -1232
$ws.Range("H32").Value = 18882444
$ws.Range("I32").Value = 20846094
$ws.Range("J32").Value = 31399.4
$ws.Range("K32").Value = 20846094
$ws.Range("L32").Value = 31399.4
$ws.Range("M32").Value = -20845807
$ws.Range("N32").Value = -31973.4
$ws.Range("H61").Value = 5800.4517
$ws.Range("I61").Value = 4848.381
$ws.Range("K61").Value = 4848.381
$ws.Range("M61").Value = -4636.381
$ws.Range("H63").Value = 3166.25
$ws.Range("I63").Value = 2599.5
$ws.Range("J63").Value = 6000
$ws.Range("K63").Value = 2599.5
$ws.Range("L63").Value = 6000
$ws.Range("M63").Value = -1913.5
$ws.Range("N63").Value = -7372
$ws.Range("H66").Value = 3166.25
$ws.Range("I66").Value = 2599.5
$ws.Range("J66").Value = 6000
$ws.Range("K66").Value = 12997.5
$ws.Range("L66").Value = 30000
$ws.Range("M66").Value = -9565.5
$ws.Range("N66").Value = -36864
$ws.Range("H74").Value = 1847.9
$ws.Range("I74").Value = 1379.625
$ws.Range("J74").Value = 3721
$ws.Range("K74").Value = 1379.625
$ws.Range("L74").Value = 3721
$ws.Range("M74").Value = -505.625
$ws.Range("N74").Value = -5469
$ws.Range("H77").Value = 1847.9
$ws.Range("I77").Value = 1379.625
$ws.Range("J77").Value = 3721
$ws.Range("K77").Value = 6898.125
$ws.Range("L77").Value = 18605
$ws.Range("M77").Value = -2530.125
$ws.Range("N77").Value = -27341
$ws.Range("H116").Value = 15800.8
$ws.Range("J116").Value = 1006
$ws.Range("L116").Value = 1006
$ws.Range("N116").Value = -5594
$ws.Range("H122").Value = 4400.8
$ws.Range("I122").Value = 3513.8125
$ws.Range("J122").Value = 5977.6665
$ws.Range("K122").Value = 10541.4375
$ws.Range("L122").Value = 17932.9995
$ws.Range("M122").Value = -8091.4375
$ws.Range("N122").Value = -22832.9995
$ws.Range("H132").Value = 3753.8718
$ws.Range("I132").Value = 3753.8718
$ws.Range("K132").Value = 11261.6154
$ws.Range("M132").Value = -8731.615399999999
$ws.Range("H136").Value = 5800.4517
$ws.Range("I136").Value = 4848.381
$ws.Range("K136").Value = 14545.143
$ws.Range("M136").Value = -11995.143
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 15800.8
$ws.Range("J3").Value = 1006
$ws.Range("L3").Value = 1006
$ws.Range("N3").Value = -1234
$ws.Range("H20").Value = 52717.1
$ws.Range("I20").Value = 2843.3845
$ws.Range("J20").Value = 145339.72
$ws.Range("K20").Value = 2843.3845
$ws.Range("L20").Value = 145339.72
$ws.Range("M20").Value = -2596.3845
$ws.Range("N20").Value = -145833.72
$ws.Range("H58").Value = 48371.5
$ws.Range("J58").Value = 45000
$ws.Range("L58").Value = 45000
$ws.Range("N58").Value = -45588
$ws.Range("H59").Value = 83997.5
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 83997.5
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 83997.5
$ws.Range("M59").Value = ""
$ws.Range("N59").Value = -85691.5
$ws.Range("H134").Value = 5088.4287
$ws.Range("I134").Value = 3613.7778
$ws.Range("J134").Value = 7742.8
$ws.Range("K134").Value = 10841.3334
$ws.Range("L134").Value = 23228.4
$ws.Range("M134").Value = -8306.3334
$ws.Range("N134").Value = -28298.4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H111").Value = 70000
$ws.Range("J111").Value = 70000
$ws.Range("L111").Value = 70000
$ws.Range("N111").Value = -78180
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = ""
$ws.Range("H122").Value = 7373921.5
$ws.Range("I122").Value = 17689356
$ws.Range("J122").Value = 5753.4287
$ws.Range("K122").Value = 53068068
$ws.Range("L122").Value = 17260.2861
$ws.Range("M122").Value = -53065618
$ws.Range("N122").Value = -22160.2861
$ws.Range("H132").Value = 3580.8635
$ws.Range("I132").Value = 3580.8635
$ws.Range("K132").Value = 10742.5905
$ws.Range("M132").Value = -8212.5905
$ws.Range("H134").Value = 3483.9583
$ws.Range("I134").Value = 1978.6111
$ws.Range("K134").Value = 5935.8333
$ws.Range("M134").Value = -3400.8333
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 172.42857
$ws.Range("I2").Value = 96
$ws.Range("J2").Value = 241.90909
$ws.Range("K2").Value = 576
$ws.Range("L2").Value = 1451.45454
$ws.Range("M2").Value = -463
$ws.Range("N2").Value = -1677.45454
$ws.Range("H37").Value = 273304.25
$ws.Range("J37").Value = 273304.25
$ws.Range("L37").Value = 819912.75
$ws.Range("N37").Value = -820136.75
$ws.Range("H45").Value = 3355
$ws.Range("I45").Value = 2220
$ws.Range("J45").Value = 4490
$ws.Range("K45").Value = 6660
$ws.Range("L45").Value = 13470
$ws.Range("M45").Value = -6128
$ws.Range("N45").Value = -14534
$ws.Range("H114").Value = 1177.9231
$ws.Range("I114").Value = 615
$ws.Range("J114").Value = 1428.1111
$ws.Range("K114").Value = 1845
$ws.Range("L114").Value = 4284.3333
$ws.Range("M114").Value = 1409
$ws.Range("N114").Value = -10792.3333
$ws.Range("H131").Value = 5716.5557
$ws.Range("I131").Value = 3931.125
$ws.Range("K131").Value = 11793.375
$ws.Range("M131").Value = -6753.375
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 16725077
$ws.Range("I80").Value = 72081.69
$ws.Range("J80").Value = 83337060
$ws.Range("K80").Value = 72081.69
$ws.Range("L80").Value = 83337060
$ws.Range("M80").Value = -71083.69
$ws.Range("N80").Value = -83339056
$ws.Range("H83").Value = 16725077
$ws.Range("I83").Value = 72081.69
$ws.Range("J83").Value = 83337060
$ws.Range("K83").Value = 360408.45
$ws.Range("L83").Value = 416685300
$ws.Range("M83").Value = -355416.45
$ws.Range("N83").Value = -416695284
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = ""
$ws.Range("H113").Value = 6125.4287
$ws.Range("I113").Value = 3813
$ws.Range("K113").Value = 3813
$ws.Range("M113").Value = -1643
$ws.Range("H132").Value = 4568.65
$ws.Range("I132").Value = 3607.9678
$ws.Range("K132").Value = 10823.9034
$ws.Range("M132").Value = -8293.903399999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 38495
$ws.Range("J54").Value = 38495
$ws.Range("L54").Value = 38495
$ws.Range("N54").Value = -39783
$ws.Range("H101").Value = 61142.57
$ws.Range("I101").Value = 48999
$ws.Range("K101").Value = 48999
$ws.Range("M101").Value = -45754
$ws.Range("H132").Value = 4741.057
$ws.Range("I132").Value = 4703.8823
$ws.Range("K132").Value = 14111.6469
$ws.Range("M132").Value = -11581.6469
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3425.7144
$ws.Range("I132").Value = 1808.8
$ws.Range("J132").Value = 4895.636
$ws.Range("K132").Value = 5426.4
$ws.Range("L132").Value = 14686.908
$ws.Range("M132").Value = -2896.4
$ws.Range("N132").Value = -19746.908
$ws.Range("H136").Value = 3702.0527
$ws.Range("I136").Value = 2925.276
$ws.Range("K136").Value = 8775.828
$ws.Range("M136").Value = -6225.828

Write-Output "Applied 211 cell updates across 8 sheets"